$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5450
$ws1.Range("F5").Value = 308
$ws1.Range("F8").Value = 349

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5450
$ws4.Range("F5").Value = 308
$ws4.Range("F9").Value = 349
